$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7, column B: change from "Dispense" (shared string) to a numeric price
# value of 1, formatted as currency with no decimals (matches existing style
# of B3, a "$"#,##0 custom format) -> displays as $1.
$ws.Range("B7").Value = 1
$ws.Range("B7").NumberFormat = $ws.Range("B3").NumberFormat

# New row 9: A9 = 6, B9 = "Dispense" (the text that used to be in B7)
$ws.Range("A9").Value = 6
$ws.Range("B9").Value = "Dispense"

# Update the active selection to B8, as reflected in the saved sheet view.
$ws.Range("B8").Select()
